$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer RestApi")

# Row 11 (RestApi #9): 비교 페이지 접속 stays, URI updated from graph -> comparison
$ws.Range("E11").Value = "/comparison/index"

# Row 12 (RestApi #10): 그래프 아이템 검색 -> 비교 아이템 검색, URI updated
$ws.Range("E12").Value = "/comparison/items"

# Row 13 (RestApi #11): 그래프 아이템 확정 -> 비교 아이템 확정, URI updated, note added
$ws.Range("E13").Value = "/comparison/item/:id"
$ws.Range("H13").Value = "프론트단에서 처리하도록 변경"

$ws.Range("C12").Value = "비교 아이템 검색"
$ws.Range("C13").Value = "비교 아이템 확정"

# Column H widens (bestFit) to accommodate the longer note text
$ws.Columns.Item(8).ColumnWidth = 28.285714285714285

$ws.Range("F14").Select()

$wb.Save()
